# Add two new columns, I (I0) and J (IF), to the worksheet.
# I0 is always 1 for every data row; IF mirrors the existing IP column (H).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers, styled the same way as the existing header cells (e.g. H1 / IP).
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Determine last used row based on column A (data rows 2..30).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $hValue = [double]($ws.Cells.Item($r, 8).Value2)  # column H (IP)
    $ws.Cells.Item($r, 9).Value2 = 1          # column I (I0)
    $ws.Cells.Item($r, 10).Value2 = $hValue   # column J (IF)
}
